# Auto-generated: apply the cryptos-list refresh described by the commit
# "Updated cryptos list on Mon Aug 28 18:36:13 UTC 2023 with GitHub Actions".
# All target cells are plain text cells (inlineStr) in the source workbook, so
# every write below keeps the cell as text -- values that look like plain
# numbers get a leading apostrophe (Excel's text-qualifier) so they are not
# silently reinterpreted as numbers, and the cell Style is reset to "Normal"
# right after so no stray number-format style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.231.78'
$ws.Range("D3").Value = '1.655.94'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("D5").Value = '''219.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").Value = '''0.5242'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.24%  '
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("D8").Value = '''0.2664'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").Value = '''0.06357'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = '''20.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("D11").Value = '''0.07711'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").Value = '''4.600'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.87%  '
$ws.Range("D13").Value = '1.621.04'
$ws.Range("E13").Value = '  -2.72%  '
$ws.Range("D14").Value = '1.884.13'
$ws.Range("D15").Value = '''0.5631'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = '0.0₅8246'
$ws.Range("E16").Value = '  +1.02%  '
$ws.Range("D17").Value = '''65.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.54%  '
$ws.Range("D18").Value = '26.232.99'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = '''4.696'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("D21").Value = '''10.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").Value = '''192.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.19%  '
$ws.Range("D23").Value = '''6.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("D25").Value = '''143.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = '  -1.16%  '
$ws.Range("D27").Value = '''7.273'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("D29").Value = '''1.513'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("D30").Value = '''0.05629'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.33%  '
$ws.Range("D31").Value = '''1.278'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").Value = '''3.508'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("D33").Value = '''3.358'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.85%  '
$ws.Range("D34").Value = '''1.585'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("D35").Value = '''0.9544'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.53%  '
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").Value = '''2.412'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.00%  '
$ws.Range("E38").Value = '  -0.85%  '
$ws.Range("D39").Value = '''0.01600'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").Value = '''6.002'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.98%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '''1.004'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.59%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''0.8425'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("D43").Value = '''101.97'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("D44").Value = '1.008.41'
$ws.Range("E44").Value = '  -6.26%  '
$ws.Range("D45").Value = '1.794.90'
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").Value = '''58.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D48").Value = '''0.05350'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.85%  '
$ws.Range("D49").Value = '0.0₈103'
$ws.Range("E49").Value = '  -2.79%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''8.037'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.4349'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.40%  '
